$d = $word.ActiveDocument

# Mapping of old text -> new text for this edit.
$replacements = @(
    @("2023-09-20 Wednesday", "2023-09-21 Thursday"),
    @("27×18=", "71×49="),
    @("77×36=", "53×74="),
    @("32×47=", "88×54="),
    @("23×90=", "21×32="),
    @("82×57=", "94×11="),
    @("56×79=", "47×59="),
    @("66×51=", "28×16="),
    @("44×41=", "55×57="),
    @("31×21=", "23×96="),
    @("60×14=", "42×66="),
    @("29×78=", "39×13="),
    @("86×84=", "55×32="),
    @("20×27=", "65×16="),
    @("34×54=", "76×56="),
    @("67×47=", "26×17="),
    @("97×35=", "32×14="),
    @("65×77=", "87×65="),
    @("42×89=", "50×82="),
    @("31×15=", "83×30="),
    @("13×15=", "48×21="),
    @("88×23=", "13×59="),
    @("33×43=", "92×54="),
    @("70×69=", "54×58="),
    @("20×24=", "74×16="),
    @("15×96=", "21×16=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
